$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 146-147; this pushes the existing rows 146-159
# down to 148-161 (preserving all of their data/formatting).
$ws.Rows("146:147").Insert()

# --- New row 146 (week of 2021-11-16 / serial 44516) ---
$ws.Range("A146").Value = 4
$ws.Range("B146").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C146").Value = "Los Lagos"
$ws.Range("D146").Value = 44516
$ws.Range("D146").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E146").Value = 10
$ws.Range("F146").Value = "Fruta"
$ws.Range("G146").Value = 100101
$ws.Range("H146").Value = "Berries"
$ws.Range("I146").Value = 100101007
$ws.Range("J146").Value = "Kiwi"
$ws.Range("K146").Value = "Hayward"
$ws.Range("L146").Value = "Especial"
$ws.Range("M146").Value = 250
$ws.Range("N146").Value = 22000
$ws.Range("O146").Value = 22000
$ws.Range("P146").Value = 22000
$ws.Range("Q146").Value = "$/caja 15 kilos"
$ws.Range("R146").Value = "Región de O'Higgins"
$ws.Range("S146").Value = 1467
$ws.Range("T146").Value = 15

# --- New row 147 (week of 2021-11-16 / serial 44516) ---
$ws.Range("A147").Value = 4
$ws.Range("B147").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C147").Value = "Los Lagos"
$ws.Range("D147").Value = 44516
$ws.Range("D147").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E147").Value = 10
$ws.Range("F147").Value = "Fruta"
$ws.Range("G147").Value = 100101
$ws.Range("H147").Value = "Berries"
$ws.Range("I147").Value = 100101007
$ws.Range("J147").Value = "Kiwi"
$ws.Range("K147").Value = "Hayward"
$ws.Range("L147").Value = "Primera"
$ws.Range("M147").Value = 600
$ws.Range("N147").Value = 15000
$ws.Range("O147").Value = 16000
$ws.Range("P147").Value = 15500
$ws.Range("Q147").Value = "$/caja 15 kilos"
$ws.Range("R147").Value = "Región de O'Higgins"
$ws.Range("S147").Value = 1033
$ws.Range("T147").Value = 15

# The row that used to be 152 (Primera, Provincia de Curicó, 2021-08-13)
# is now row 154 after the insert; its quality label is corrected from
# "Primera" to "Segunda".
$ws.Range("L154").Value = "Segunda"
